# Applies the "Creado grafico de tipos de modelo" edit:
# Inserts a new "MAE" column (D) before the existing "Tipo" column,
# shifting "Tipo" to column E, and refreshes the MSE/R2 values (and
# adds the new MAE values) for all 13 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts the old D ("Tipo") to E.
$ws.Columns.Item(4).Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "MAE"

# Row data: Enfermedad, MSE (B), R2 (C), MAE (D)
$data = @(
    @{ Row = 2;  B = 0.4858260388868722;  C = 0.9903261718552575;  D = 0.5688422646190253 },
    @{ Row = 3;  B = 0.2357639951764815;  C = 0.995391986785256;   D = 0.3843755001218409 },
    @{ Row = 4;  B = 0.354228512740891;   C = 0.9931854310920872;  D = 0.4698397028756816 },
    @{ Row = 5;  B = 0.406772042895757;   C = 0.9919790122346771;  D = 0.4952415476002354 },
    @{ Row = 6;  B = 0.7384315664895728;  C = 0.9783103305233276;  D = 0.6639375116161906 },
    @{ Row = 7;  B = 0.3631266388613507;  C = 0.9949814915291805;  D = 0.4965276629928901 },
    @{ Row = 8;  B = 0.1556483113569549;  C = 0.9983957660431915;  D = 0.3309224161203504 },
    @{ Row = 9;  B = 0.4564187238483384;  C = 0.9972812438886709;  D = 0.5581454810666486 },
    @{ Row = 10; B = 0.1058018669828489;  C = 0.9980672504508001;  D = 0.2369360161650128 },
    @{ Row = 11; B = 0.1986423075790961;  C = 0.9853235826675191;  D = 0.3425956061944328 },
    @{ Row = 12; B = 0.05081863269278668; C = 0.9985055107115712;  D = 0.1661962461621603 },
    @{ Row = 13; B = 0.09974499197663596; C = 0.9990538740505828;  D = 0.2214877650927295 },
    @{ Row = 14; B = 0.08996850467433244; C = 0.9987818819825056;  D = 0.2397878346241068 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
}
